# This script reproduces the commit "Elimina antiguos EC y agrega nuevos y
# modifica Antigua BD": it appends a new "Estado de Cuenta" data row for
# period 2509 (copying the format/values of the previous last row, 2508),
# restores the previous last row to the regular (non-closing) row style, and
# updates the two summary cells (total "Valor Mora" and "Cant. Periodos").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row right after the current last data row (row 21),
#    which pushes the closing rows ("___", "NOMBRE DEL REPRESENTANTE LEGAL",
#    "FIRMA DEL REPRESENTANTE LEGAL") down from rows 26/27 to rows 27/28.
$ws.Rows("22:22").Insert()

# 2. Copy row 21 (values + full formatting, including the special "closing"
#    bottom border used for the last row of the table) into the newly
#    inserted row 22, so row 22 becomes the new last row of the table.
$ws.Range("B21:J21").Copy($ws.Range("B22:J22"))

# 3. Copy row 20's formatting + values (a normal "middle of table" row) onto
#    row 21, so row 21 loses the special closing border that used to belong
#    to the old last row.
$ws.Range("B20:J20").Copy($ws.Range("B21:J21"))

# 4. Step 3 overwrote row 21's "Periodo Mora" value with row 20's (2507), so
#    restore it back to its real value, 2508 (stored as text, like the rest
#    of the period column).
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2508"

# 5. Set the new row's "Periodo Mora" value to the new period, 2509.
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2509"

# 6. Update the summary fields: total "Valor Mora" grew from 309374 to
#    366314, and "Cant. Periodos" grew from 6 to 7 (one new period added).
$ws.Range("E11").Value = 366314
$ws.Range("F13").Value = 7

$excel.CutCopyMode = $false
